$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2:A3").Value = "Pfizer_mono_Janssen_Novavax"
$ws.Range("A4:A5").Value = "Janssen_Novavax"
$ws.Range("A6").Value = "Pfizer_bi_Moderna_mono_Moderna_bi_Novavax"
$ws.Range("A7:A8").Value = "Pfizer_bi_Moderna_bi_Novavax"
$ws.Range("A9").Value = "Moderna_bi_Novavax"
$ws.Range("A10").Value = "Pfizer_bi_Novavax"
$ws.Range("A11:A16").Value = "Pfizer_mono_Novavax"
$ws.Range("A17:A109").Value = "Novavax"
$ws.Range("A110").Value = "Pfizer_mono_Pfizer_bi_Moderna_bi_Janssen"
$ws.Range("A111").Value = "Pfizer_bi_Moderna_bi_Janssen"
$ws.Range("A112").Value = "Pfizer_mono_Moderna_bi_Janssen"
$ws.Range("A113").Value = "Pfizer_mono_Moderna_mono_Janssen"
$ws.Range("A114:A115").Value = "Pfizer_mono_Pfizer_bi_Janssen"
$ws.Range("A116:A119").Value = "Pfizer_bi_Janssen"
$ws.Range("A120:A129").Value = "Pfizer_mono_Janssen"
$ws.Range("A130:A160").Value = "Janssen"
$ws.Range("A161:A166").Value = "Pfizer_mono_Pfizer_bi_Moderna_bi"
$ws.Range("A167:A186").Value = "Pfizer_bi_Moderna_bi"
$ws.Range("A187").Value = "Pfizer_mono_Moderna_bi"
$ws.Range("A188").Value = "Moderna_bi"
$ws.Range("A189").Value = "Pfizer_mono_Moderna_mono"
$ws.Range("A190:A200").Value = "Moderna_mono"
$ws.Range("A201:A202").Value = "Pfizer_mono_Pfizer_bi"
$ws.Range("A203:A214").Value = "Pfizer_bi"
$ws.Range("A215:A229").Value = "Pfizer_mono"
